# "Update Work Week and Social Spending"
#
# The Clio-Infra "GDP per Capita" extract for Equatorial Guinea (country code
# 226) was refreshed with a newer revision of the source series:
#   - every existing yearly value (1950-2008) is replaced by the new figure
#   - eight more years (2009-2016) are appended to the Data sheet
# The Metadata sheet content (citations/links) is unchanged.

$wb     = $excel.ActiveWorkbook
$wsData = $wb.Worksheets.Item("Data")

# GDP per Capita, in order, for years 1950..2016 (one entry per year).
$gdpValues = @(
    "797", "827", "845", "867", "901", "915", "940", "958",
    "974", "1015", "1101", "1191", "1310", "1486", "1690", "1898",
    "1953", "2058", "2102", "2015", "1934", "1989", "2007", "2204",
    "2466", "3080", "3453", "3668", "3451", "2882", "2613", "2603",
    "2533", "2534", "2437", "2662", "2534", "2579", "2577", "2479",
    "2474", "2427.88024951205", "2690.63309412509", "2880.20400818885",
    "3081.15250420458", "3364.63258044284", "4426.33782610114", "8445.64499672705",
    "9712.22621077315", "11690.1351585363", "12914.8512921307", "21167.7413539634",
    "24845.0422260709", "27571.1590436026", "33419.8857944737", "35569.7076236273",
    "36615.7057037176", "41277.7556536689", "47562.3195289018", "47161.2847380668",
    "42037.8974036049", "43841", "46255", "43209",
    "41907", "37822", "33317"
)

$firstYear = 1950
$firstDataRow = 2      # row 2 holds year 1950
$lastExistingRow = 60  # row 60 holds year 2008 (last row present before the edit)
$lastRow = $firstDataRow + $gdpValues.Length - 1   # row 68 holds year 2016

for ($row = $firstDataRow; $row -le $lastRow; $row++) {
    $year = $firstYear + ($row - $firstDataRow)
    $value = $gdpValues[$row - $firstDataRow]

    if ($row -gt $lastExistingRow) {
        # Brand-new row for a year that didn't exist before (2009-2016).
        $wsData.Range("A$row").Value = 226
        $wsData.Range("B$row").Value = "Equatorial Guinea"
        $wsData.Range("C$row").Value = "GDP per Capita"
        $wsData.Range("D$row").Value = $year
    }

    # The GDP figures are stored as text (not numbers) in this workbook, same
    # as the original data. A leading apostrophe forces Excel to keep a
    # numeric-looking entry as text; ClearFormats() then drops the transient
    # "quote prefix" cell format so the cell keeps the sheet's default style.
    $wsData.Range("E$row").Value = "'" + $value
    $wsData.Range("E$row").ClearFormats()
}
